$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: In the "HBR[%] = 100 * Sum_{j=1..M} x_beneficial,j" formula,
# split the nary subscript run "beneficial,j" into two runs:
# "harmful" and ",j" (so the text reads x_harmful,j).
# ---------------------------------------------------------------------
$targetHBR = [string][char]119867 + [char]119861 + [char]119877   # math-italic "HBR"
$target1 = $null
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $t = $om.Range.Text
    if ($t -like "$targetHBR*") {
        $target1 = $om
    }
}
if ($target1 -eq $null) {
    throw "Could not locate the HBR[%] summation OMath formula"
}

$xml1 = '<m:oMathPara xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<m:oMath>' +
        '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>HBR[%]= 100</m:t></m:r>' +
        '<m:nary>' +
        '<m:naryPr>' +
        '<m:chr m:val="&#8721;"/>' +
        '<m:limLoc m:val="undOvr"/>' +
        '<m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr>' +
        '</m:naryPr>' +
        '<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>j=1</m:t></m:r></m:sub>' +
        '<m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>M</m:t></m:r></m:sup>' +
        '<m:e>' +
        '<m:sSub>' +
        '<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>' +
        '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e>' +
        '<m:sub>' +
        '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>harmful</m:t></m:r>' +
        '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>,j</m:t></m:r>' +
        '</m:sub>' +
        '</m:sSub>' +
        '</m:e>' +
        '</m:nary>' +
        '</m:oMath>' +
        '</m:oMathPara>'

$target1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# Edit 2: the standalone math abbreviation run that precedes the
# "Harmful bacteria ratio" legend text currently (mis)reads "BBR"; it
# must read "HBR", split into two runs: "H" and "BR".
# ---------------------------------------------------------------------
$targetBBR = [string][char]119861 + [char]119861 + [char]119877   # math-italic "BBR"
$target2 = $null
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $t = $om.Range.Text
    if ($t -like "$targetBBR*" -and $t.Length -le 10) {
        $full = $d.Range($om.Range.Start, $om.Range.End).Text
        if ($full -like "*Harmful*") {
            $target2 = $om
        }
    }
}
if ($target2 -eq $null) {
    throw "Could not locate the standalone BBR OMath run preceding 'Harmful bacteria ratio'"
}

$xml2 = '<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>H</m:t></m:r>' +
        '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>BR</m:t></m:r>' +
        '</m:oMath>'

$target2.Range.InsertXML($xml2)
